$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount cell (T2) with the new value from the daily upload.
$ws.Range("T2").Value = 131518
